$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto price/volume refresh (GitHub Actions bot)
# D column (Price) holds plain text values (not locale numbers), so force
# text formatting before/after the write to avoid Excel auto-converting
# numeric-looking strings (e.g. "159.62") into real numbers.

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '69.387.63'
$cell.Style = "Normal"
$ws.Range('E2').Value = '  +0.07%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '3.674.94'
$cell.Style = "Normal"
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  +0.00%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '641.32'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -5.41%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '159.62'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('E7').Value = '  +0.02%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.497'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('E10').Value = '  -1.30%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.447'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('E12').Value = '  +0.14%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '4.295.72'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  -0.28%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '32.62'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  +0.89%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '3.678.22'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -0.54%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '69.373.43'
$cell.Style = "Normal"
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('E17').Value = '  -0.01%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '15.97'
$cell.Style = "Normal"
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('E19').Value = '  +0.23%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '465.97'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -0.55%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '9.92'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  +0.50%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '0.648'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -0.68%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '79.30'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -0.86%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '3.821.93'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  -0.31%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  +2.62%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '10.89'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -0.21%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '9.06'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('E32').Value = '  -0.17%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '26.87'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('E34').Value = '  +3.53%  '
$ws.Range('E35').Value = '  -1.72%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '3.668.68'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('E37').Value = '  +1.44%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('E39').Value = '  -6.56%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '177.88'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +4.36%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('E44').Value = '  -1.81%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '46.76'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('E46').Value = '  +1.34%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '27.32'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -2.96%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '0.000269'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -3.00%  '
$ws.Range('E49').Value = '  -3.43%  '
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('E51').Value = '  -3.40%  '
